$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '95.125.75'
$ws.Cells.Item(2, 5).Value = '  -0.93%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.543.74'
$ws.Cells.Item(3, 5).Value = '  -0.51%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.01%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''235.21'
$ws.Cells.Item(5, 5).Value = '  -1.82%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''648.51'
$ws.Cells.Item(6, 5).Value = '  +1.70%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''1.46'
$ws.Cells.Item(7, 5).Value = '  -1.49%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -1.37%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.06%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''0.989'
$ws.Cells.Item(10, 5).Value = '  -3.27%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '3.539.99'
$ws.Cells.Item(11, 5).Value = '  -0.50%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +0.04%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''41.95'
$ws.Cells.Item(13, 5).Value = '  -3.31%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''6.48'
$ws.Cells.Item(14, 5).Value = '  +1.30%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '4.204.49'
$ws.Cells.Item(15, 5).Value = '  -0.56%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '95.011.04'
$ws.Cells.Item(16, 5).Value = '  -1.04%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '''0.0000251'
$ws.Cells.Item(17, 5).Value = '  -1.18%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.552.38'
$ws.Cells.Item(18, 5).Value = '  +0.67%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -1.77%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''12.55'
$ws.Cells.Item(20, 5).Value = '  -2.68%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''17.66'
$ws.Cells.Item(21, 5).Value = '  -1.67%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +0.13%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''502.17'
$ws.Cells.Item(23, 5).Value = '  -2.58%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''0.470'
$ws.Cells.Item(24, 5).Value = '  -6.77%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -1.39%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''6.62'
$ws.Cells.Item(26, 5).Value = '  -2.21%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''94.77'
$ws.Cells.Item(27, 5).Value = '  -2.04%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '3.734.67'
$ws.Cells.Item(28, 5).Value = '  -0.50%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '''12.32'
$ws.Cells.Item(29, 5).Value = '  -0.49%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''3.00'
$ws.Cells.Item(30, 5).Value = '  -2.09%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'Hedera'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(31, 4).Value = '''0.140'
$ws.Cells.Item(31, 5).Value = '  -3.95%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(32, 4).Value = '''11.27'
$ws.Cells.Item(32, 5).Value = '  -2.28%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'Dai'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(33, 4).Value = '''0.998'
$ws.Cells.Item(33, 5).Value = '  -0.15%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''0.999'
$ws.Cells.Item(34, 5).Value = '  -0.01%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -3.87%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''31.65'
$ws.Cells.Item(36, 5).Value = '  +4.94%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -2.58%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''8.27'
$ws.Cells.Item(38, 5).Value = '  +6.87%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'Fetch.AI'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(39, 4).Value = '''1.56'
$ws.Cells.Item(39, 5).Value = '  +7.41%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'Bittensor'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(40, 4).Value = '''564.57'
$ws.Cells.Item(40, 5).Value = '  -3.12%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.03%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -1.62%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''0.892'
$ws.Cells.Item(43, 5).Value = '  -4.25%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''1.73'
$ws.Cells.Item(44, 5).Value = '  -0.66%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''2.29'
$ws.Cells.Item(45, 5).Value = '  +4.89%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(46, 4).Value = '''23.40'
$ws.Cells.Item(46, 5).Value = '  -2.05%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'EnergySwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(47, 4).Value = '''33.58'
$ws.Cells.Item(47, 5).Value = '  +31.64%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''5.56'
$ws.Cells.Item(48, 5).Value = '  -0.98%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'VeChain'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(49, 4).Value = '''0.0408'
$ws.Cells.Item(49, 5).Value = '  -4.54%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'MantraDAO'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Cells.Item(50, 4).Value = '''3.54'
$ws.Cells.Item(50, 5).Value = '  -0.55%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -1.78%  '
